# Generate Report for Handback
# The fc32db19-... file has now been handed back (in all languages), so:
#   - the Overview sheet's status for that file flips from "Ready for
#     handoff" to "Handed back: in sync with en-US"
#   - each language sheet's row for that file gets its Status updated the
#     same way, plus the "Latest Target File" / "Latest Handback File"
#     hyperlinks populated (they mirror the Source/Handoff file links) and
#     the "Latest Handback DateTime" stamped with the real handback time
#     (replacing the 0001-01-01 00:00:00 placeholder).

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: columns are File Name | zh-cn | de-de
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B5").Value = $statusHandedBack
$overview.Range("C5").Value = $statusHandedBack

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B5").Value = $statusHandedBack

$zhcn.Hyperlinks.Add(
    $zhcn.Cells.Item(5, 5),
    "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/4471afbd8f869e1e97eef4b7646283df1cfeb180/e2e/fc32db19-d4c5-4ee4-b250-092a647c2d39.md",
    "",
    "",
    "fc32db19-d4c5-4ee4-b250-092a647c2d39.md"
) | Out-Null

$zhcn.Hyperlinks.Add(
    $zhcn.Cells.Item(5, 6),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/4471afbd8f869e1e97eef4b7646283df1cfeb180/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/fc32db19-d4c5-4ee4-b250-092a647c2d39.4471afbd8f869e1e97eef4b7646283df1cfeb180.zh-cn.xlf",
    "",
    "",
    "fc32db19-d4c5-4ee4-b250-092a647c2d39.4471afbd8f869e1e97eef4b7646283df1cfeb180.zh-cn.xlf"
) | Out-Null

$zhcn.Range("G5").Value = "2016-02-22 04:51:57"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B5").Value = $statusHandedBack

$dede.Hyperlinks.Add(
    $dede.Cells.Item(5, 5),
    "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/4471afbd8f869e1e97eef4b7646283df1cfeb180/e2e/fc32db19-d4c5-4ee4-b250-092a647c2d39.md",
    "",
    "",
    "fc32db19-d4c5-4ee4-b250-092a647c2d39.md"
) | Out-Null

$dede.Hyperlinks.Add(
    $dede.Cells.Item(5, 6),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/4471afbd8f869e1e97eef4b7646283df1cfeb180/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/fc32db19-d4c5-4ee4-b250-092a647c2d39.4471afbd8f869e1e97eef4b7646283df1cfeb180.de-de.xlf",
    "",
    "",
    "fc32db19-d4c5-4ee4-b250-092a647c2d39.4471afbd8f869e1e97eef4b7646283df1cfeb180.de-de.xlf"
) | Out-Null

$dede.Range("G5").Value = "2016-02-22 04:52:21"

Write-Host "Handback report generated."
